# Auto-generated Excel COM-interop edit script
# Applies the cryptos.xlsx "Updated symbol list" diff (Mon Jan 23 15:56:50 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "304.13"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.06%"
$ws.Range("E2").NumberFormat = "General"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "35.61"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.04%"
$ws.Range("E3").NumberFormat = "General"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.051"
$ws.Range("D4").NumberFormat = "General"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07869"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "0.42%"
$ws.Range("E5").NumberFormat = "General"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.132"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-3.26%"
$ws.Range("E6").NumberFormat = "General"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.133"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "2.77%"
$ws.Range("E7").NumberFormat = "General"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "7.940"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.92%"
$ws.Range("E8").NumberFormat = "General"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9238"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.65%"
$ws.Range("E9").NumberFormat = "General"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09838"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.40%"
$ws.Range("E10").NumberFormat = "General"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1850"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-2.29%"
$ws.Range("E11").NumberFormat = "General"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08587"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.53%"
$ws.Range("E12").NumberFormat = "General"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03582"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.13%"
$ws.Range("E13").NumberFormat = "General"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09932"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.16%"
$ws.Range("E14").NumberFormat = "General"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001435"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-3.56%"
$ws.Range("E15").NumberFormat = "General"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005720"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.32%"
$ws.Range("E16").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.16%"
$ws.Range("E17").NumberFormat = "General"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.753"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "15.98%"
$ws.Range("E18").NumberFormat = "General"
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3374"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-1.08%"
$ws.Range("E19").NumberFormat = "General"
$ws.Range("B20").Value = "ProBitToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1347"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.87%"
$ws.Range("E20").NumberFormat = "General"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.176"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "8.56%"
$ws.Range("E21").NumberFormat = "General"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2208"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "0.04%"
$ws.Range("E22").NumberFormat = "General"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04577"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.72%"
$ws.Range("E23").NumberFormat = "General"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.76%"
$ws.Range("E24").NumberFormat = "General"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004832"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "7.59%"
$ws.Range("E25").NumberFormat = "General"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001302"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-7.22%"
$ws.Range("E26").NumberFormat = "General"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0004751"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "74.78%"
$ws.Range("E27").NumberFormat = "General"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01856"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.07%"
$ws.Range("E39").NumberFormat = "General"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04719"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-1.39%"
$ws.Range("E40").NumberFormat = "General"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007793"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-3.14%"
$ws.Range("E41").NumberFormat = "General"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1390"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.99%"
$ws.Range("E42").NumberFormat = "General"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.007737"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "3.02%"
$ws.Range("E43").NumberFormat = "General"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002224"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "4.32%"
$ws.Range("E44").NumberFormat = "General"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.01134"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "8.93%"
$ws.Range("E45").NumberFormat = "General"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006303"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.51%"
$ws.Range("E46").NumberFormat = "General"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("E47").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.15%"
$ws.Range("E48").NumberFormat = "General"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "51.80"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "37.06%"
$ws.Range("E49").NumberFormat = "General"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.001901"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-29.34%"
$ws.Range("E50").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
$ws.Range("E51").NumberFormat = "General"

Write-Output "Applied 102 cell updates"
